{"js": "// Replace the 100 arithmetic answers in the (20 rows x 5 cols) table with a\n// new set of values. Each table cell holds a single paragraph/run, so we\n// target each cell's first paragraph range directly and use\n// Word.InsertLocation.replace \u2014 this swaps only the `<w:t>` text and keeps\n// the existing run/paragraph formatting (font, size, justification) intact.\n\nconst after = [\n  [\"51-50=1\", \"5+35=40\", \"25-21=4\", \"25-25=0\", \"84-4=80\"],\n  [\"56-28=28\", \"97-41=56\", \"39-22=17\", \"48-24=24\", \"94-74=20\"],\n  [\"7+75=82\", \"89-70=19\", \"95-82=13\", \"86-51=35\", \"90-11=79\"],\n  [\"80-68=12\", \"34+64=98\", \"32+11=43\", \"66-5=61\", \"61+9=70\"],\n  [\"36-24=12\", \"7+83=90\", \"29-20=9\", \"39-14=25\", \"25+32=57\"],\n  [\"66-8=58\", \"39+47=86\", \"19+6=25\", \"81-62=19\", \"44+28=72\"],\n  [\"9+80=89\", \"44-8=36\", \"91-13=78\", \"14+14=28\", \"36+40=76\"],\n  [\"36-14=22\", \"77-58=19\", \"96-60=36\", \"23-3=20\", \"16-14=2\"],\n  [\"27+39=66\", \"63-22=41\", \"36+33=69\", \"95-71=24\", \"85+0=85\"],\n  [\"75-71=4\", \"26+43=69\", \"39-33=6\", \"82-21=61\", \"46+1=47\"],\n  [\"34+40=74\", \"70+4=74\", \"43+37=80\", \"34+52=86\", \"24+48=72\"],\n  [\"79-44=35\", \"58-40=18\", \"58-27=31\", \"91-36=55\", \"91-25=66\"],\n  [\"99-14=85\", \"26+33=59\", \"85-46=39\", \"88-70=18\", \"13+73=86\"],\n  [\"54+18=72\", \"52+18=70\", \"64-33=31\", \"10+64=74\", \"29+46=75\"],\n  [\"40+14=54\", \"60-8=52\", \"23-16=7\", \"44-36=8\", \"53-22=31\"],\n  [\"7+58=65\", \"9+19=28\", \"25+23=48\", \"61-21=40\", \"15-8=7\"],\n  [\"77-2=75\", \"46-8=38\", \"50+18=68\", \"47+4=51\", \"61+33=94\"],\n  [\"40+22=62\", \"93-11=82\", \"14+28=42\", \"38+33=71\", \"56-19=37\"],\n  [\"35+64=99\", \"0+1=1\", \"34-13=21\", \"82-44=38\", \"34+30=64\"],\n  [\"21+35=56\", \"96-13=83\", \"23+70=93\", \"79-6=73\", \"4-3=1\"],\n];\n\nconst table = context.document.body.tables.getFirst();\n\nfor (let r = 0; r < after.length; r++) {\n  const rowValues = after[r];\n  for (let c = 0; c < rowValues.length; c++) {\n    const cell = table.getCell(r, c);\n    cell.body.paragraphs.load(\"items\");\n    await context.sync();\n\n    const paragraph = cell.body.paragraphs.items[0];\n    const range = paragraph.getRange();\n    range.insertText(rowValues[c], Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 100 arithmetic answers in the (20 rows x 5 cols) table with\n# a new set of values, preserving each cell's existing run/paragraph\n# formatting (font, size, justification) by assigning straight to\n# Cell.Range.Text (this keeps the surrounding w:rPr/w:pPr untouched).\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n$values = @(\n    @(\"51-50=1\", \"5+35=40\", \"25-21=4\", \"25-25=0\", \"84-4=80\"),\n    @(\"56-28=28\", \"97-41=56\", \"39-22=17\", \"48-24=24\", \"94-74=20\"),\n    @(\"7+75=82\", \"89-70=19\", \"95-82=13\", \"86-51=35\", \"90-11=79\"),\n    @(\"80-68=12\", \"34+64=98\", \"32+11=43\", \"66-5=61\", \"61+9=70\"),\n    @(\"36-24=12\", \"7+83=90\", \"29-20=9\", \"39-14=25\", \"25+32=57\"),\n    @(\"66-8=58\", \"39+47=86\", \"19+6=25\", \"81-62=19\", \"44+28=72\"),\n    @(\"9+80=89\", \"44-8=36\", \"91-13=78\", \"14+14=28\", \"36+40=76\"),\n    @(\"36-14=22\", \"77-58=19\", \"96-60=36\", \"23-3=20\", \"16-14=2\"),\n    @(\"27+39=66\", \"63-22=41\", \"36+33=69\", \"95-71=24\", \"85+0=85\"),\n    @(\"75-71=4\", \"26+43=69\", \"39-33=6\", \"82-21=61\", \"46+1=47\"),\n    @(\"34+40=74\", \"70+4=74\", \"43+37=80\", \"34+52=86\", \"24+48=72\"),\n    @(\"79-44=35\", \"58-40=18\", \"58-27=31\", \"91-36=55\", \"91-25=66\"),\n    @(\"99-14=85\", \"26+33=59\", \"85-46=39\", \"88-70=18\", \"13+73=86\"),\n    @(\"54+18=72\", \"52+18=70\", \"64-33=31\", \"10+64=74\", \"29+46=75\"),\n    @(\"40+14=54\", \"60-8=52\", \"23-16=7\", \"44-36=8\", \"53-22=31\"),\n    @(\"7+58=65\", \"9+19=28\", \"25+23=48\", \"61-21=40\", \"15-8=7\"),\n    @(\"77-2=75\", \"46-8=38\", \"50+18=68\", \"47+4=51\", \"61+33=94\"),\n    @(\"40+22=62\", \"93-11=82\", \"14+28=42\", \"38+33=71\", \"56-19=37\"),\n    @(\"35+64=99\", \"0+1=1\", \"34-13=21\", \"82-44=38\", \"34+30=64\"),\n    @(\"21+35=56\", \"96-13=83\", \"23+70=93\", \"79-6=73\", \"4-3=1\")\n)\n\nfor ($r = 1; $r -le $values.Length; $r++) {\n    $row = $values[$r - 1]\n    for ($c = 1; $c -le $row.Length; $c++) {\n        $tbl.Cell($r, $c).Range.Text = $row[$c - 1]\n    }\n}\n"}
